# Fruta / hortaliza, semanal
# Update weekly price records for Terminal Hortofrutícola Agro Chillán - Frambuesa
# by rewriting the data rows (2-5) with the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44995
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 5500
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 5750
$ws.Range("S2").Value = 2875

# Row 3
$ws.Range("D3").Value = 44991
$ws.Range("M3").Value = 50

# Row 4
$ws.Range("D4").Value = 45008
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("S4").Value = 3500

# Row 5
$ws.Range("D5").Value = 45008
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000
